$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 48.484375
$ws.Range("H2").Value = 145.453125
$ws.Range("I2").Value = 0.7776469276297807
$ws.Range("J2").Value = 0.7776469276297806
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 75.59011833333334
$ws.Range("N2").Value = 226.770355
$ws.Range("O2").Value = 0.6588374259037486
$ws.Range("P2").Value = 0.6588374259037486
$ws.Range("Q2").Value = 3664.939643567709
$ws.Range("R2").Value = 32984.45679210938
$ws.Range("S2").Value = 0.5123429000615634
$ws.Range("T2").Value = 0.5123429000615634
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 48.484375
$ws.Range("H3").Value = 145.453125
$ws.Range("I3").Value = 0.7776469276297807
$ws.Range("J3").Value = 0.7776469276297806
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 12.15310033333333
$ws.Range("N3").Value = 36.459301
$ws.Range("O3").Value = 0.1059254505338229
$ws.Range("P3").Value = 0.1059254505338229
$ws.Range("Q3").Value = 589.2354739739583
$ws.Range("R3").Value = 5303.119265765625
$ws.Range("S3").Value = 0.08237260116542772
$ws.Range("T3").Value = 0.0823726011654277
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 48.484375
$ws.Range("H4").Value = 145.453125
$ws.Range("I4").Value = 0.7776469276297807
$ws.Range("J4").Value = 0.7776469276297806
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1727356666666667
$ws.Range("N4").Value = 0.5182070000000001
$ws.Range("O4").Value = 0.001505550255743542
$ws.Range("P4").Value = 0.001505550255743542
$ws.Range("Q4").Value = 8.374980838541669
$ws.Range("R4").Value = 75.37482754687501
$ws.Range("S4").Value = 0.001170786530771196
$ws.Range("T4").Value = 0.001170786530771196
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 48.484375
$ws.Range("H5").Value = 145.453125
$ws.Range("I5").Value = 0.7776469276297807
$ws.Range("J5").Value = 0.7776469276297806
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 26.81662666666667
$ws.Range("N5").Value = 80.44988000000001
$ws.Range("O5").Value = 0.233731573306685
$ws.Range("P5").Value = 0.233731573306685
$ws.Range("Q5").Value = 1300.187383541667
$ws.Range("R5").Value = 11701.686451875
$ws.Range("S5").Value = 0.1817606398720184
$ws.Range("T5").Value = 0.1817606398720184
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 11.47148733333333
$ws.Range("H6").Value = 34.414462
$ws.Range("I6").Value = 0.1839926137051496
$ws.Range("J6").Value = 0.1839926137051496
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 75.59011833333334
$ws.Range("N6").Value = 226.770355
$ws.Range("O6").Value = 0.6588374259037486
$ws.Range("P6").Value = 0.6588374259037486
$ws.Range("Q6").Value = 867.1310849860012
$ws.Range("R6").Value = 7804.17976487401
$ws.Range("S6").Value = 0.1212212199988035
$ws.Range("T6").Value = 0.1212212199988035
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 11.47148733333333
$ws.Range("H7").Value = 34.414462
$ws.Range("I7").Value = 0.1839926137051496
$ws.Range("J7").Value = 0.1839926137051496
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 12.15310033333333
$ws.Range("N7").Value = 36.459301
$ws.Range("O7").Value = 0.1059254505338229
$ws.Range("P7").Value = 0.1059254505338229
$ws.Range("Q7").Value = 139.4141365345625
$ws.Range("R7").Value = 1254.727228811062
$ws.Range("S7").Value = 0.01948950050161362
$ws.Range("T7").Value = 0.01948950050161361
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 11.47148733333333
$ws.Range("H8").Value = 34.414462
$ws.Range("I8").Value = 0.1839926137051496
$ws.Range("J8").Value = 0.1839926137051496
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.1727356666666667
$ws.Range("N8").Value = 0.5182070000000001
$ws.Range("O8").Value = 0.001505550255743542
$ws.Range("P8").Value = 0.001505550255743542
$ws.Range("Q8").Value = 1.981535012181556
$ws.Range("R8").Value = 17.833815109634
$ws.Range("S8").Value = 0.0002770101266187108
$ws.Range("T8").Value = 0.0002770101266187107
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 11.47148733333333
$ws.Range("H9").Value = 34.414462
$ws.Range("I9").Value = 0.1839926137051496
$ws.Range("J9").Value = 0.1839926137051496
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 26.81662666666667
$ws.Range("N9").Value = 80.44988000000001
$ws.Range("O9").Value = 0.233731573306685
$ws.Range("P9").Value = 0.233731573306685
$ws.Range("Q9").Value = 307.6265931293956
$ws.Range("R9").Value = 2768.63933816456
$ws.Range("S9").Value = 0.04300488307811374
$ws.Range("T9").Value = 0.04300488307811374
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.356432
$ws.Range("H10").Value = 1.069296
$ws.Range("I10").Value = 0.005716857228930723
$ws.Range("J10").Value = 0.005716857228930722
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 75.59011833333334
$ws.Range("N10").Value = 226.770355
$ws.Range("O10").Value = 0.6588374259037486
$ws.Range("P10").Value = 0.6588374259037486
$ws.Range("Q10").Value = 26.94273705778667
$ws.Range("R10").Value = 242.48463352008
$ws.Range("S10").Value = 0.003766479500967955
$ws.Range("T10").Value = 0.003766479500967954
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.356432
$ws.Range("H11").Value = 1.069296
$ws.Range("I11").Value = 0.005716857228930723
$ws.Range("J11").Value = 0.005716857228930722
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 12.15310033333333
$ws.Range("N11").Value = 36.459301
$ws.Range("O11").Value = 0.1059254505338229
$ws.Range("P11").Value = 0.1059254505338229
$ws.Range("Q11").Value = 4.331753858010667
$ws.Range("R11").Value = 38.985784722096
$ws.Range("S11").Value = 0.0006055606776120293
$ws.Range("T11").Value = 0.0006055606776120292
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.356432
$ws.Range("H12").Value = 1.069296
$ws.Range("I12").Value = 0.005716857228930723
$ws.Range("J12").Value = 0.005716857228930722
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.1727356666666667
$ws.Range("N12").Value = 0.5182070000000001
$ws.Range("O12").Value = 0.001505550255743542
$ws.Range("P12").Value = 0.001505550255743542
$ws.Range("Q12").Value = 0.06156851914133335
$ws.Range("R12").Value = 0.5541166722720001
$ws.Range("S12").Value = [double]"8.60701586306597E-06"
$ws.Range("T12").Value = [double]"8.607015863065968E-06"
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.356432
$ws.Range("H13").Value = 1.069296
$ws.Range("I13").Value = 0.005716857228930723
$ws.Range("J13").Value = 0.005716857228930722
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 26.81662666666667
$ws.Range("N13").Value = 80.44988000000001
$ws.Range("O13").Value = 0.233731573306685
$ws.Range("P13").Value = 0.233731573306685
$ws.Range("Q13").Value = 9.558303876053335
$ws.Range("R13").Value = 86.02473488448001
$ws.Range("S13").Value = 0.001336210034487673
$ws.Range("T13").Value = 0.001336210034487673
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 2.035248333333333
$ws.Range("H14").Value = 6.105745
$ws.Range("I14").Value = 0.03264360143613892
$ws.Range("J14").Value = 0.03264360143613892
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 75.59011833333334
$ws.Range("N14").Value = 226.770355
$ws.Range("O14").Value = 0.6588374259037486
$ws.Range("P14").Value = 0.6588374259037486
$ws.Range("Q14").Value = 153.8446623543861
$ws.Range("R14").Value = 1384.601961189475
$ws.Range("S14").Value = 0.02150682634241368
$ws.Range("T14").Value = 0.02150682634241368
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 2.035248333333333
$ws.Range("H15").Value = 6.105745
$ws.Range("I15").Value = 0.03264360143613892
$ws.Range("J15").Value = 0.03264360143613892
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 12.15310033333333
$ws.Range("N15").Value = 36.459301
$ws.Range("O15").Value = 0.1059254505338229
$ws.Range("P15").Value = 0.1059254505338229
$ws.Range("Q15").Value = 24.73457719824944
$ws.Range("R15").Value = 222.611194784245
$ws.Range("S15").Value = 0.003457788189169565
$ws.Range("T15").Value = 0.003457788189169564
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 2.035248333333333
$ws.Range("H16").Value = 6.105745
$ws.Range("I16").Value = 0.03264360143613892
$ws.Range("J16").Value = 0.03264360143613892
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.1727356666666667
$ws.Range("N16").Value = 0.5182070000000001
$ws.Range("O16").Value = 0.001505550255743542
$ws.Range("P16").Value = 0.001505550255743542
$ws.Range("Q16").Value = 0.3515599776905556
$ws.Range("R16").Value = 3.164039799215
$ws.Range("S16").Value = [double]"4.914658249056923E-05"
$ws.Range("T16").Value = [double]"4.914658249056923E-05"
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 2.035248333333333
$ws.Range("H17").Value = 6.105745
$ws.Range("I17").Value = 0.03264360143613892
$ws.Range("J17").Value = 0.03264360143613892
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 26.81662666666667
$ws.Range("N17").Value = 80.44988000000001
$ws.Range("O17").Value = 0.233731573306685
$ws.Range("P17").Value = 0.233731573306685
$ws.Range("Q17").Value = 54.57849472895555
$ws.Range("R17").Value = 491.2064525606
$ws.Range("S17").Value = 0.007629840322065111
$ws.Range("T17").Value = 0.007629840322065111
